$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the reference list: drop stray hyphens that were being used in
# place of spaces/colons within several citation titles (column A).

$ws.Range("A16").Value = "Hargittai E 2002 - Second Level Digital Divide: Differences in People's Online Skills"
$ws.Range("A23").Value = "Selwyn N 2006-1- Digital division or digital decision? A study of non users and low users of computers"
$ws.Range("A24").Value = "Chinn M 2007 - The determinants of the global digital divide: a cross country analysis of computer and internet penetration"
$ws.Range("A25").Value = "Kiiski S 2002 - Cross country diffusion of the Internet"
$ws.Range("A27").Value = "Corrocher N 2002 - Measuring the Digital Divide: A Framework for the Analysis of Cross Country Differences"
$ws.Range("A28").Value = "Chen W 2004 - The Global Digital Divide: Within and Between Countries"
$ws.Range("A31").Value = "Beilock R 2003 - An Exploratory Model of Inter Country Internet Diffusion"

# Reflect the reviewer's working selection (the whole cleaned table body) and
# scroll the sheet back to the top of the data.
$ws.Range("A2:E32").Select()
